$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1833952912019826
$ws.Range("C2").Value = 0.5712515489467163
$ws.Range("J2").Value = 0.02478314745972739
$ws.Range("P2").Value = 0.1412639405204461
$ws.Range("S2").Value = 0.07930607187112763
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.03885480572597137
$ws.Range("J3").Value = 0.03476482617586912
$ws.Range("P3").Value = 0.7464212678936605
$ws.Range("S3").Value = 0.1738241308793456
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6615384615384615
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.0502092050209205
$ws.Range("D6").Value = 0.01115760111576011
$ws.Range("F6").Value = 0.06276150627615062
$ws.Range("J6").Value = 0.2259414225941423
$ws.Range("O6").Value = 0.03207810320781032
$ws.Range("Q6").Value = 0.1645746164574617
$ws.Range("R6").Value = 0.08507670850767085
$ws.Range("S6").Value = 0.3682008368200837
$ws.Range("B7").Value = 0.09027777777777778
$ws.Range("D7").Value = 0.015625
$ws.Range("E7").Value = 0.003472222222222222
$ws.Range("F7").Value = 0.0642361111111111
$ws.Range("J7").Value = 0.1319444444444444
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.1736111111111111
$ws.Range("R7").Value = 0.08159722222222222
$ws.Range("S7").Value = 0.4079861111111111
$ws.Range("B8").Value = 0.08157524613220815
$ws.Range("D8").Value = 0.01687763713080169
$ws.Range("F8").Value = 0.06610407876230662
$ws.Range("J8").Value = 0.129395218002813
$ws.Range("O8").Value = 0.02250351617440225
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.1075949367088608
$ws.Range("S8").Value = 0.4092827004219409
$ws.Range("B9").Value = 0.084375
$ws.Range("D9").Value = 0.0140625
$ws.Range("E9").Value = 0.0015625
$ws.Range("F9").Value = 0.0625
$ws.Range("J9").Value = 0.1140625
$ws.Range("O9").Value = 0.0234375
$ws.Range("Q9").Value = 0.16875
$ws.Range("R9").Value = 0.1015625
$ws.Range("S9").Value = 0.4296875
$ws.Range("B10").Value = 0.0927734375
$ws.Range("D10").Value = 0.01953125
$ws.Range("E10").Value = 0.001220703125
$ws.Range("F10").Value = 0.072509765625
$ws.Range("J10").Value = 0.131591796875
$ws.Range("O10").Value = 0.01806640625
$ws.Range("Q10").Value = 0.227783203125
$ws.Range("R10").Value = 0.088623046875
$ws.Range("S10").Value = 0.347900390625
$ws.Range("G11").Value = 0.1427003293084523
$ws.Range("J11").Value = 0.09440175631174534
$ws.Range("K11").Value = 0.1931942919868277
$ws.Range("L11").Value = 0.5609220636663008
$ws.Range("S11").Value = 0.008781558726673985
$ws.Range("G12").Value = 0.7421150278293135
$ws.Range("J12").Value = 0.1855287569573284
$ws.Range("K12").Value = 0.007421150278293136
$ws.Range("L12").Value = 0.04823747680890538
$ws.Range("S12").Value = 0.01669758812615955
$ws.Range("G13").Value = 0.6176470588235294
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.04901960784313725
$ws.Range("F15").Value = 0.01517241379310345
$ws.Range("H15").Value = 0.1544827586206897
$ws.Range("I15").Value = 0.05793103448275862
$ws.Range("J15").Value = 0.383448275862069
$ws.Range("K15").Value = 0.05517241379310345
$ws.Range("M15").Value = 0.008275862068965517
$ws.Range("O15").Value = 0.06758620689655172
$ws.Range("S15").Value = 0.2579310344827586
$ws.Range("F16").Value = 0.01848428835489834
$ws.Range("H16").Value = 0.1719038817005545
$ws.Range("I16").Value = 0.09611829944547134
$ws.Range("J16").Value = 0.4288354898336414
$ws.Range("K16").Value = 0.09426987060998152
$ws.Range("M16").Value = 0.01848428835489834
$ws.Range("N16").Value = 0.001848428835489834
$ws.Range("O16").Value = 0.04990757855822551
$ws.Range("S16").Value = 0.1201478743068392
$ws.Range("F17").Value = 0.01343183344526528
$ws.Range("H17").Value = 0.1732706514439221
$ws.Range("I17").Value = 0.1007387508394896
$ws.Range("J17").Value = 0.4392209536601746
$ws.Range("K17").Value = 0.09603760913364674
$ws.Range("M17").Value = 0.01276024177300202
$ws.Range("N17").Value = 0.002014775016789792
$ws.Range("O17").Value = 0.06111484217595702
$ws.Range("S17").Value = 0.1014103425117528
$ws.Range("F18").Value = 0.0249266862170088
$ws.Range("H18").Value = 0.1891495601173021
$ws.Range("I18").Value = 0.08064516129032258
$ws.Range("J18").Value = 0.4178885630498534
$ws.Range("K18").Value = 0.09824046920821114
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("N18").Value = 0.001466275659824047
$ws.Range("O18").Value = 0.0747800586510264
$ws.Range("S18").Value = 0.0967741935483871
$ws.Range("F19").Value = 0.01471747700394218
$ws.Range("H19").Value = 0.2149802890932983
$ws.Range("I19").Value = 0.08988173455978975
$ws.Range("J19").Value = 0.3779237844940867
$ws.Range("K19").Value = 0.1114323258869908
$ws.Range("M19").Value = 0.01576872536136662
$ws.Range("N19").Value = 0.0005256241787122207
$ws.Range("O19").Value = 0.07279894875164257
$ws.Range("S19").Value = 0.1019710906701708
